# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.046.42'
$ws.Range("E2").Value = '  -1.18%  '
$ws.Range("D3").Value = '1.793.40'
$ws.Range("E3").Value = '  -0.10%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = "'316.78"
$ws.Range("E5").Value = '  +0.79%  '
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = "'0.5343"
$ws.Range("E7").Value = '  -1.80%  '
$ws.Range("D8").Value = "'0.3763"
$ws.Range("E8").Value = '  -1.42%  '
$ws.Range("E9").Value = '  -1.83%  '
$ws.Range("D10").Value = "'42.06"
$ws.Range("E10").Value = '  -2.23%  '
$ws.Range("E11").Value = '  -2.77%  '
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = '  -0.09%  '
$ws.Range("E13").Value = '  -2.62%  '
$ws.Range("D14").Value = "'6.122"
$ws.Range("E14").Value = '  -1.05%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '1.798.57'
$ws.Range("E15").Value = '  +0.16%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = "'7.236"
$ws.Range("E16").Value = '  -2.08%  '
$ws.Range("D17").Value = "'89.10"
$ws.Range("E17").Value = '  -2.56%  '
$ws.Range("D18").Value = "'0.00001058"
$ws.Range("E18").Value = '  -1.18%  '
$ws.Range("D19").Value = "'0.06494"
$ws.Range("E19").Value = '  +0.71%  '
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("E21").Value = '  -0.32%  '
$ws.Range("D22").Value = "'5.898"
$ws.Range("E22").Value = '  -0.92%  '
$ws.Range("D23").Value = '28.069.15'
$ws.Range("E23").Value = '  -1.17%  '
$ws.Range("D24").Value = "'11.15"
$ws.Range("E24").Value = '  -2.51%  '
$ws.Range("E25").Value = '  -1.76%  '
$ws.Range("E26").Value = '  -2.43%  '
$ws.Range("D27").Value = "'20.29"
$ws.Range("E27").Value = '  -1.91%  '
$ws.Range("D28").Value = '1.996.42'
$ws.Range("E28").Value = '  -0.30%  '
$ws.Range("D29").Value = "'2.311"
$ws.Range("E29").Value = '  -3.12%  '
$ws.Range("D30").Value = "'121.22"
$ws.Range("E30").Value = '  -1.82%  '
$ws.Range("D31").Value = "'1.120"
$ws.Range("E31").Value = '  -0.48%  '
$ws.Range("D32").Value = "'0.1064"
$ws.Range("E32").Value = '  +4.31%  '
$ws.Range("D33").Value = "'3.666"
$ws.Range("E33").Value = '  -0.45%  '
$ws.Range("D34").Value = "'5.564"
$ws.Range("E34").Value = '  -3.28%  '
$ws.Range("D35").Value = "'0.2248"
$ws.Range("E35").Value = '  -4.13%  '
$ws.Range("E36").Value = '  -5.15%  '
$ws.Range("D37").Value = "'0.02288"
$ws.Range("E37").Value = '  -1.44%  '
$ws.Range("D38").Value = "'5.022"
$ws.Range("E38").Value = '  -2.51%  '
$ws.Range("D39").Value = "'8.473"
$ws.Range("E39").Value = '  -3.10%  '
$ws.Range("D40").Value = "'0.6186"
$ws.Range("E40").Value = '  -3.36%  '
$ws.Range("D41").Value = "'11.16"
$ws.Range("E41").Value = '  -4.26%  '
$ws.Range("D42").Value = "'1.448"
$ws.Range("E42").Value = '  +1.84%  '
$ws.Range("D43").Value = "'1.178"
$ws.Range("E43").Value = '  +1.93%  '
$ws.Range("D44").Value = "'13.31"
$ws.Range("E44").Value = '  -2.32%  '
$ws.Range("D45").Value = "'3.670"
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").Value = "'0.5786"
$ws.Range("E46").Value = '  -3.34%  '
$ws.Range("D47").Value = "'125.09"
$ws.Range("E47").Value = '  -1.13%  '
$ws.Range("D48").Value = "'1.191"
$ws.Range("E48").Value = '  +3.42%  '
$ws.Range("E49").Value = '  -3.81%  '
$ws.Range("E50").Value = '  -2.11%  '
$ws.Range("D51").Value = "'71.36"
$ws.Range("E51").Value = '  -2.16%  '

